# TP05 - Form données de projet
# "ajouts de la phase TEST1 dans le form de données de projet"
#
# Fill in the "Projet 5 / Réel" (column K) actual-effort figures for each
# development phase (rows 13-21). All the downstream totals, ratios and
# percentages (T13:T24, U13:U24, K24, K53/K54/K56/K61-K67/K70/K75-K80, ...)
# are formulas and recompute automatically on recalculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K13").Value = 60
$ws.Range("K14").Value = 30
$ws.Range("K15").Value = 30
$ws.Range("K16").Value = 30
$ws.Range("K17").Value = 75
$ws.Range("K18").Value = 15
$ws.Range("K19").Value = 45
$ws.Range("K20").Value = 12
$ws.Range("K21").Value = 8

# Restore the cursor/selection to where the author left it.
$ws.Range("L20").Select()
